$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 2 currently has no row record at all (data jumps from row 1 to
#        row 3). Touch a harmless row-level property (value unchanged) so the
#        engine emits a bare <row r="2"/> element, matching the target file. ---
$ws.Rows.Item(2).OutlineLevel = 0

# --- 2. Row 154 currently ends with two explicit-but-empty cells (H154, I154).
#        Drop them completely so the row ends at G154. ---
$ws.Range("H154:I154").ClearContents()

# --- 3. Append four new incident rows (155-158) below the existing data. ---

# Column A holds plain date-as-text values ("2024-05-22"); a direct .Value
# assignment would be auto-parsed into a real date serial number by Excel,
# which is not what the source data is (plain text). Route it through a
# text formula and paste-special-values so it lands as literal text without
# picking up a date number format / style.
function Set-TextDate($cellRef, $text) {
    $ws.Range($cellRef).Formula = '="' + $text + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

Set-TextDate "A155" "2024-05-22"
$ws.Range("B155").Value = "10:31:45"
$ws.Range("C155").Value = "-"
$ws.Range("D155").Value = "-"
$ws.Range("E155").Value = "-"
$ws.Range("F155").Value = "NOK Soldadura metal"
$ws.Range("G155").Value = "-"
$ws.Range("H155").Value = "10:31:48"
$ws.Range("I155").Value = "0:00:03"

Set-TextDate "A156" "2024-05-22"
$ws.Range("B156").Value = "10:38:28"
$ws.Range("C156").Value = "-"
$ws.Range("D156").Value = "-"
$ws.Range("E156").Value = "-"
$ws.Range("F156").Value = "-"
$ws.Range("G156").Value = "Colisión placas"
$ws.Range("H156").Value = "10:38:30"
$ws.Range("I156").Value = "0:00:02"

Set-TextDate "A157" "2024-05-22"
$ws.Range("B157").Value = "10:38:44"
$ws.Range("C157").Value = "-"
$ws.Range("D157").Value = "-"
$ws.Range("E157").Value = "-"
$ws.Range("F157").Value = "-"
$ws.Range("G157").Value = "Fallo dispensación glue"
$ws.Range("H157").Value = "10:38:47"
$ws.Range("I157").Value = "0:00:03"

Set-TextDate "A158" "2024-05-22"
$ws.Range("B158").Value = "10:38:49"
$ws.Range("C158").Value = "-"
$ws.Range("D158").Value = "-"
$ws.Range("E158").Value = "-"
$ws.Range("F158").Value = "-"
$ws.Range("G158").Value = "Error en sensor de salida"
$ws.Range("H158").Value = "10:39:53"
$ws.Range("I158").Value = "0:01:04"

Write-Host "edit complete"
